# "add tabel format baru" - renumber the table labels (Tabel 4.2.3/4.2.4/4.2.5
# -> Tabel 4.2.5/4.2.6/4.2.7) and roll the reporting year forward from 2020 to
# 2021 across the four side-by-side "Bab 4" sub-tables, then reset the sheet
# view back to the top-left cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 2 (H1/I1/I2) -------------------------------------------------
# H1 held the plain label "Tabel 4.2.3" -> becomes "Tabel 4.2.5"
$ws.Range("H1").Value = "Tabel 4.2.5"

# Indonesian + English titles: bump the year 2020 -> 2021
$ws.Range("I1").Value = "Banyaknya Bayi yang Diimunisasi Menurut Jenis dan Desa/Kelurahan di Kecamatan Samaturu, 2021"
$ws.Range("I2").Value = "Number of Immunized Babies by Types of Immunization and Kelurahan/Village Samaturu Subdistrict, 2021"

# --- Table 1 (A1/B1/B2) --------------------------------------------------
$ws.Range("B1").Value = "Jumlah Tenaga Kesehatan Menurut Kelurahan/Desa in Kecamatan Samaturu. 2021"
$ws.Range("B2").Value = "Number of Medical Personnel by Kelurahan/ Village in Samaturu Subdistrict, 2021"

# --- Table 3 (P1/Q1/Q2) ---------------------------------------------------
# P1 is rich text: "Tabel" (inherits the underlined cell font) + " 4.2.4."
# (explicit, non-underlined run) -> keep the same split, just bump the
# table number to " 4.2.6."
$ws.Range("P1").Value = "Tabel 4.2.6."
$p1run = $ws.Range("P1").Characters(6, 7)
$p1run.Font.Underline = 0

$ws.Range("Q1").Value = "Banyaknya Ibu Melahirkan dan Kelahiran Ditolong Tenaga Kesehatan Menurut Desa/Kelurahan di Kecamatan Samaturu, 2021"
$ws.Range("Q2").Value = "Number of Woman Giving Brth and Birth Assisted by Paramedics by Kelurahan/Village in Samaturu Subdistrict, 2021"

# --- Table 4 (W1/X1/X2) ---------------------------------------------------
# W1 is rich text just like P1, " 4.2.5." -> " 4.2.7."
$ws.Range("W1").Value = "Tabel 4.2.7."
$w1run = $ws.Range("W1").Characters(6, 7)
$w1run.Font.Underline = 0

$ws.Range("X1").Value = "Banyaknya Pasangan Usia Subur dan Peserta KB Menurut Desa/Kelurahan di Kecamatan Samaturu, 2021"
$ws.Range("X2").Value = "Number of Fertile Age Couples and Family Planning Members by Kelurahan/Village in Samaturu Subdistrict, 2021"

# --- Reset the sheet view: scroll back to the top-left and select A1 -----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
